$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$sh = $s.Shapes.Item("PoljeZBesedilom 7")
$sh.Width = 121.8898
$sh.Height = 153.0709

$sh = $s.Shapes.Item("PoljeZBesedilom 1")
$sh.Width = 121.8898
$sh.Height = 153.0709

$sh = $s.Shapes.Item("PoljeZBesedilom 16")
$sh.Width = 121.8898
$sh.Height = 153.0709

$sh = $s.Shapes.Item("PoljeZBesedilom 20")
$sh.Width = 121.8898
$sh.Height = 153.0709

$sh = $s.Shapes.Item("PoljeZBesedilom 26")
$sh.Width = 121.8898
$sh.Height = 153.0709

$sh = $s.Shapes.Item("PoljeZBesedilom 34")
$sh.Width = 121.8898
$sh.Height = 153.0709

$sh = $s.Shapes.Item("PoljeZBesedilom 36")
$sh.Width = 121.8898
$sh.Height = 153.0709

$sh = $s.Shapes.Item("PoljeZBesedilom 44")
$sh.Width = 121.8898
$sh.Height = 153.0709

$sh = $s.Shapes.Item("Pravokotnik 3")
$sh.Left = 408.2028
$sh.Top = 231.68551
$sh.Width = 31.181102
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 40")
$sh.Left = 325.528
$sh.Top = 231.68551
$sh.Width = 76.535435
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 41")
$sh.Left = 258.02614
$sh.Top = 231.68551
$sh.Width = 31.181102
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 45")
$sh.Left = 175.35134125
$sh.Top = 231.68551
$sh.Width = 76.535435
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 46")
$sh.Left = 105.7727
$sh.Top = 231.68551
$sh.Width = 31.181102
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 47")
$sh.Left = 23.24701
$sh.Top = 231.68551
$sh.Width = 76.535435
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 48")
$sh.Left = 105.7727
$sh.Top = 475.93945
$sh.Width = 31.181102
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 51")
$sh.Left = 23.24701
$sh.Top = 475.93945
$sh.Width = 76.535435
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 52")
$sh.Left = 258.02614
$sh.Top = 475.93945
$sh.Width = 31.181102
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 53")
$sh.Left = 174.598504
$sh.Top = 475.93945
$sh.Width = 76.535435
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 54")
$sh.Left = 408.350632
$sh.Top = 475.93945
$sh.Width = 31.181102
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 55")
$sh.Left = 325.082535
$sh.Top = 475.93945
$sh.Width = 76.535435
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 56")
$sh.Left = 557.9245
$sh.Top = 475.93945
$sh.Width = 31.181102
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 57")
$sh.Left = 475.0255
$sh.Top = 475.93945
$sh.Width = 76.535435
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 58")
$sh.Left = 705.9328
$sh.Top = 475.93945
$sh.Width = 31.181102
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 59")
$sh.Left = 623.0338
$sh.Top = 475.93945
$sh.Width = 76.535435
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 60")
$sh.Left = 705.9328
$sh.Top = 231.68551
$sh.Width = 31.181102
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 61")
$sh.Left = 623.258
$sh.Top = 231.68551
$sh.Width = 76.535435
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 62")
$sh.Left = 558.12064
$sh.Top = 231.68551
$sh.Width = 31.181102
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 63")
$sh.Left = 476.4571
$sh.Top = 231.68551
$sh.Width = 76.535435
$sh.Height = 14.17323

$sh = $s.Shapes.Item("Pravokotnik 63")
$tr = $sh.TextFrame.TextRange
$c = $tr.Characters(1, 9)
$c.Text = "GENERATE "
